$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated stats) ---
$ws.Range("G9").Value = 0.0364278880542198
$ws.Range("L9").Value = 0.00295
$ws.Range("G10").Value = 0.0364278880542198
$ws.Range("L10").Value = 0.00295
$ws.Range("G13").Value = 0.695540677966102
$ws.Range("I13").Value = 1.34107
$ws.Range("L13").Value = 0.5128
$ws.Range("G14").Value = 0.695540677966102
$ws.Range("I14").Value = 1.34107
$ws.Range("L14").Value = 0.5128
$ws.Range("G26").Value = 0.0236776937096546
$ws.Range("L26").Value = 0.00295
$ws.Range("G27").Value = 0.0236776937096546
$ws.Range("L27").Value = 0.00295
$ws.Range("F30").Value = 0.69285
$ws.Range("G30").Value = 0.728998333333333
$ws.Range("I30").Value = 1.3387
$ws.Range("L30").Value = 0.70885
$ws.Range("F31").Value = 0.69285
$ws.Range("G31").Value = 0.728998333333333
$ws.Range("I31").Value = 1.3387
$ws.Range("L31").Value = 0.70885
$ws.Range("G39").Value = 1354.50933358794
$ws.Range("H39").Value = 11633.5600152764
$ws.Range("G40").Value = 1354.50933358794
$ws.Range("H40").Value = 11633.5600152764
$ws.Range("G41").Value = 1354.50933358794
$ws.Range("H41").Value = 11633.5600152764
$ws.Range("G42").Value = 1354.50933358794
$ws.Range("H42").Value = 11633.5600152764
$ws.Range("G43").Value = 0.0285840987201545
$ws.Range("G44").Value = 0.0285840987201545
$ws.Range("G47").Value = 0.7792416666666671
$ws.Range("I47").Value = 1.3387
$ws.Range("G48").Value = 0.7792416666666671
$ws.Range("I48").Value = 1.3387
$ws.Range("G56").Value = 1844.40933358794
$ws.Range("G57").Value = 1844.40933358794
$ws.Range("G58").Value = 1844.40933358794
$ws.Range("G59").Value = 1844.40933358794
$ws.Range("G60").Value = 0.0403080000294403
$ws.Range("G61").Value = 0.0403080000294403
$ws.Range("G62").Value = 0.660680283928732
$ws.Range("G63").Value = 0.660680283928732
$ws.Range("G64").Value = 0.769526420024975
$ws.Range("G65").Value = 0.769526420024975
$ws.Range("G73").Value = 1717.19593246231
$ws.Range("G74").Value = 1717.19593246231
$ws.Range("G75").Value = 1717.19593246231
$ws.Range("G76").Value = 1717.19593246231
$ws.Range("G77").Value = 0.0461581491493698
$ws.Range("G78").Value = 0.0461581491493698
$ws.Range("G79").Value = 0.623552233061927
$ws.Range("G80").Value = 0.623552233061927
$ws.Range("G81").Value = 0.725359071211838
$ws.Range("G82").Value = 0.725359071211838
$ws.Range("G90").Value = 1885.90620715994
$ws.Range("G91").Value = 1885.90620715994
$ws.Range("G92").Value = 1885.90620715994
$ws.Range("G93").Value = 1885.90620715994
$ws.Range("G94").Value = 0.0578616257086667
$ws.Range("G95").Value = 0.0578616257086667
$ws.Range("G96").Value = 0.534443511343945
$ws.Range("L96").Value = 0.11389
$ws.Range("G97").Value = 0.534443511343945
$ws.Range("L97").Value = 0.11389
$ws.Range("G98").Value = 0.638254349512572
$ws.Range("G99").Value = 0.638254349512572
$ws.Range("G107").Value = 1893.54413819442
$ws.Range("G108").Value = 1893.54413819442
$ws.Range("G109").Value = 1893.54413819442
$ws.Range("G110").Value = 1893.54413819442
$ws.Range("G111").Value = 0.0583099860549522
$ws.Range("G112").Value = 0.0583099860549522
$ws.Range("G113").Value = 0.54190614929734
$ws.Range("L113").Value = 0.11389
$ws.Range("G114").Value = 0.54190614929734
$ws.Range("L114").Value = 0.11389
$ws.Range("G115").Value = 0.637323315029814
$ws.Range("G116").Value = 0.637323315029814

# --- Append new rows 121-137 (2019 - 2023 reporting period) ---
# Row 121
$ws.Range("A121").Value = "Oroua at Mangawhata"
$ws.Range("B121").Value = "Visual Clarity (Sediment class 1)"
$ws.Range("C121").Value = "D"
$ws.Range("D121").Value = "2019 - 2023"
$ws.Range("E121").Value = "RepSite"
$ws.Range("F121").Value = 0.1
$ws.Range("G121").Value = 0.18530303030303
$ws.Range("H121").Value = 0.49
$ws.Range("I121").Value = 0.4885
$ws.Range("L121").Value = 0.12
$ws.Range("M121").Value = 0.39
$ws.Range("N121").Value = 0.4772
$ws.Range("O121").Value = 1806649.9
$ws.Range("P121").Value = 5525867
$ws.Range("Q121").Value = "Manawatu District"
$ws.Range("R121").Value = "Manawatū"
$ws.Range("S121").Value = "Oroua"
$ws.Range("T121").Value = "Mana_12c"
$ws.Range("U121").Value = "m"

# Row 122
$ws.Range("A122").Value = "Oroua at Mangawhata"
$ws.Range("B122").Value = "DRP (95th Percentile)"
$ws.Range("C122").Value = "C"
$ws.Range("D122").Value = "2019 - 2023"
$ws.Range("E122").Value = "RepSite"
$ws.Range("F122").Value = 0.019
$ws.Range("G122").Value = 0.0222413793103448
$ws.Range("H122").Value = 0.117
$ws.Range("I122").Value = 0.0428
$ws.Range("L122").Value = 0.019
$ws.Range("M122").Value = 0.031
$ws.Range("N122").Value = 0.03958
$ws.Range("O122").Value = 1806649.9
$ws.Range("P122").Value = 5525867
$ws.Range("Q122").Value = "Manawatu District"
$ws.Range("R122").Value = "Manawatū"
$ws.Range("S122").Value = "Oroua"
$ws.Range("T122").Value = "Mana_12c"
$ws.Range("U122").Value = "mg/L"

# Row 123
$ws.Range("A123").Value = "Oroua at Mangawhata"
$ws.Range("B123").Value = "DRP (Median)"
$ws.Range("C123").Value = "D"
$ws.Range("D123").Value = "2019 - 2023"
$ws.Range("E123").Value = "RepSite"
$ws.Range("F123").Value = 0.019
$ws.Range("G123").Value = 0.0222413793103448
$ws.Range("H123").Value = 0.117
$ws.Range("I123").Value = 0.0428
$ws.Range("L123").Value = 0.019
$ws.Range("M123").Value = 0.031
$ws.Range("N123").Value = 0.03958
$ws.Range("O123").Value = 1806649.9
$ws.Range("P123").Value = 5525867
$ws.Range("Q123").Value = "Manawatu District"
$ws.Range("R123").Value = "Manawatū"
$ws.Range("S123").Value = "Oroua"
$ws.Range("T123").Value = "Mana_12c"
$ws.Range("U123").Value = "mg/L"

# Row 124
$ws.Range("A124").Value = "Oroua at Mangawhata"
$ws.Range("B124").Value = "E coli (>260)"
$ws.Range("C124").Value = "E"
$ws.Range("D124").Value = "2019 - 2023"
$ws.Range("E124").Value = "RepSite"
$ws.Range("F124").Value = 361.5
$ws.Range("G124").Value = 1656.20529461638
$ws.Range("H124").Value = 49000
$ws.Range("I124").Value = 3768
$ws.Range("J124").Value = 36.2068965517241
$ws.Range("K124").Value = 68.9655172413793
$ws.Range("L124").Value = 556
$ws.Range("M124").Value = 1464
$ws.Range("N124").Value = 2586
$ws.Range("O124").Value = 1806649.9
$ws.Range("P124").Value = 5525867
$ws.Range("Q124").Value = "Manawatu District"
$ws.Range("R124").Value = "Manawatū"
$ws.Range("S124").Value = "Oroua"
$ws.Range("T124").Value = "Mana_12c"
$ws.Range("U124").Value = "% exceedances over 260/100 mL"

# Row 125
$ws.Range("A125").Value = "Oroua at Mangawhata"
$ws.Range("B125").Value = "E coli (>540)"
$ws.Range("C125").Value = "E"
$ws.Range("D125").Value = "2019 - 2023"
$ws.Range("E125").Value = "RepSite"
$ws.Range("F125").Value = 361.5
$ws.Range("G125").Value = 1656.20529461638
$ws.Range("H125").Value = 49000
$ws.Range("I125").Value = 3768
$ws.Range("J125").Value = 36.2068965517241
$ws.Range("K125").Value = 68.9655172413793
$ws.Range("L125").Value = 556
$ws.Range("M125").Value = 1464
$ws.Range("N125").Value = 2586
$ws.Range("O125").Value = 1806649.9
$ws.Range("P125").Value = 5525867
$ws.Range("Q125").Value = "Manawatu District"
$ws.Range("R125").Value = "Manawatū"
$ws.Range("S125").Value = "Oroua"
$ws.Range("T125").Value = "Mana_12c"
$ws.Range("U125").Value = "% exceedances over 540/100 mL"

# Row 126
$ws.Range("A126").Value = "Oroua at Mangawhata"
$ws.Range("B126").Value = "E coli (Median)"
$ws.Range("C126").Value = "E"
$ws.Range("D126").Value = "2019 - 2023"
$ws.Range("E126").Value = "RepSite"
$ws.Range("F126").Value = 361.5
$ws.Range("G126").Value = 1656.20529461638
$ws.Range("H126").Value = 49000
$ws.Range("I126").Value = 3768
$ws.Range("J126").Value = 36.2068965517241
$ws.Range("K126").Value = 68.9655172413793
$ws.Range("L126").Value = 556
$ws.Range("M126").Value = 1464
$ws.Range("N126").Value = 2586
$ws.Range("O126").Value = 1806649.9
$ws.Range("P126").Value = 5525867
$ws.Range("Q126").Value = "Manawatu District"
$ws.Range("R126").Value = "Manawatū"
$ws.Range("S126").Value = "Oroua"
$ws.Range("T126").Value = "Mana_12c"
$ws.Range("U126").Value = "E. coli/100 mL"

# Row 127
$ws.Range("A127").Value = "Oroua at Mangawhata"
$ws.Range("B127").Value = "E coli (95th Percentile)"
$ws.Range("C127").Value = "E"
$ws.Range("D127").Value = "2019 - 2023"
$ws.Range("E127").Value = "RepSite"
$ws.Range("F127").Value = 361.5
$ws.Range("G127").Value = 1656.20529461638
$ws.Range("H127").Value = 49000
$ws.Range("I127").Value = 3768
$ws.Range("J127").Value = 36.2068965517241
$ws.Range("K127").Value = 68.9655172413793
$ws.Range("L127").Value = 556
$ws.Range("M127").Value = 1464
$ws.Range("N127").Value = 2586
$ws.Range("O127").Value = 1806649.9
$ws.Range("P127").Value = 5525867
$ws.Range("Q127").Value = "Manawatu District"
$ws.Range("R127").Value = "Manawatū"
$ws.Range("S127").Value = "Oroua"
$ws.Range("T127").Value = "Mana_12c"
$ws.Range("U127").Value = "E. coli/100 mL"

# Row 128
$ws.Range("A128").Value = "Oroua at Mangawhata"
$ws.Range("B128").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C128").Value = "B"
$ws.Range("D128").Value = "2019 - 2023"
$ws.Range("E128").Value = "RepSite"
$ws.Range("F128").Value = 0.02355
$ws.Range("G128").Value = 0.0606936875191013
$ws.Range("H128").Value = 0.508041035176414
$ws.Range("I128").Value = 0.25202
$ws.Range("L128").Value = 0.00758
$ws.Range("M128").Value = 0.09322
$ws.Range("N128").Value = 0.17678
$ws.Range("O128").Value = 1806649.9
$ws.Range("P128").Value = 5525867
$ws.Range("Q128").Value = "Manawatu District"
$ws.Range("R128").Value = "Manawatū"
$ws.Range("S128").Value = "Oroua"
$ws.Range("T128").Value = "Mana_12c"
$ws.Range("U128").Value = "mg NH4-N/L"

# Row 129
$ws.Range("A129").Value = "Oroua at Mangawhata"
$ws.Range("B129").Value = "Ammoniacal-N (Median)"
$ws.Range("C129").Value = "A"
$ws.Range("D129").Value = "2019 - 2023"
$ws.Range("E129").Value = "RepSite"
$ws.Range("F129").Value = 0.02355
$ws.Range("G129").Value = 0.0606936875191013
$ws.Range("H129").Value = 0.508041035176414
$ws.Range("I129").Value = 0.25202
$ws.Range("L129").Value = 0.00758
$ws.Range("M129").Value = 0.09322
$ws.Range("N129").Value = 0.17678
$ws.Range("O129").Value = 1806649.9
$ws.Range("P129").Value = 5525867
$ws.Range("Q129").Value = "Manawatu District"
$ws.Range("R129").Value = "Manawatū"
$ws.Range("S129").Value = "Oroua"
$ws.Range("T129").Value = "Mana_12c"
$ws.Range("U129").Value = "mg NH4-N/L"

# Row 130
$ws.Range("A130").Value = "Oroua at Mangawhata"
$ws.Range("B130").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C130").Value = "A"
$ws.Range("D130").Value = "2019 - 2023"
$ws.Range("E130").Value = "RepSite"
$ws.Range("F130").Value = 0.5335
$ws.Range("G130").Value = 0.510992356193892
$ws.Range("H130").Value = 1.19
$ws.Range("I130").Value = 1.086
$ws.Range("L130").Value = 0.1115
$ws.Range("M130").Value = 0.92284
$ws.Range("N130").Value = 1.01356
$ws.Range("O130").Value = 1806649.9
$ws.Range("P130").Value = 5525867
$ws.Range("Q130").Value = "Manawatu District"
$ws.Range("R130").Value = "Manawatū"
$ws.Range("S130").Value = "Oroua"
$ws.Range("T130").Value = "Mana_12c"
$ws.Range("U130").Value = "mg NO3-N/L"

# Row 131
$ws.Range("A131").Value = "Oroua at Mangawhata"
$ws.Range("B131").Value = "Nitrate-N (Median)"
$ws.Range("C131").Value = "A"
$ws.Range("D131").Value = "2019 - 2023"
$ws.Range("E131").Value = "RepSite"
$ws.Range("F131").Value = 0.5335
$ws.Range("G131").Value = 0.510992356193892
$ws.Range("H131").Value = 1.19
$ws.Range("I131").Value = 1.086
$ws.Range("L131").Value = 0.1115
$ws.Range("M131").Value = 0.92284
$ws.Range("N131").Value = 1.01356
$ws.Range("O131").Value = 1806649.9
$ws.Range("P131").Value = 5525867
$ws.Range("Q131").Value = "Manawatu District"
$ws.Range("R131").Value = "Manawatū"
$ws.Range("S131").Value = "Oroua"
$ws.Range("T131").Value = "Mana_12c"
$ws.Range("U131").Value = "mg NO3-N/L"

# Row 132
$ws.Range("A132").Value = "Oroua at Mangawhata"
$ws.Range("B132").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("D132").Value = "2019 - 2023"
$ws.Range("E132").Value = "RepSite"
$ws.Range("F132").Value = 0.6565
$ws.Range("G132").Value = 0.60615090123671
$ws.Range("H132").Value = 1.379
$ws.Range("I132").Value = 1.2218
$ws.Range("L132").Value = 0.082
$ws.Range("M132").Value = 1.06628
$ws.Range("N132").Value = 1.14716
$ws.Range("O132").Value = 1806649.9
$ws.Range("P132").Value = 5525867
$ws.Range("Q132").Value = "Manawatu District"
$ws.Range("R132").Value = "Manawatū"
$ws.Range("S132").Value = "Oroua"
$ws.Range("T132").Value = "Mana_12c"
$ws.Range("U132").Value = "g/m3"

# Row 133
$ws.Range("A133").Value = "Oroua at Mangawhata"
$ws.Range("B133").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("D133").Value = "2019 - 2023"
$ws.Range("E133").Value = "RepSite"
$ws.Range("F133").Value = 0.6565
$ws.Range("G133").Value = 0.60615090123671
$ws.Range("H133").Value = 1.379
$ws.Range("I133").Value = 1.2218
$ws.Range("L133").Value = 0.082
$ws.Range("M133").Value = 1.06628
$ws.Range("N133").Value = 1.14716
$ws.Range("O133").Value = 1806649.9
$ws.Range("P133").Value = 5525867
$ws.Range("Q133").Value = "Manawatu District"
$ws.Range("R133").Value = "Manawatū"
$ws.Range("S133").Value = "Oroua"
$ws.Range("T133").Value = "Mana_12c"
$ws.Range("U133").Value = "g/m3"

# Row 134
$ws.Range("A134").Value = "Oroua at Mangawhata"
$ws.Range("B134").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("D134").Value = "2019 - 2023"
$ws.Range("E134").Value = "RepSite"
$ws.Range("F134").Value = 0.865
$ws.Range("G134").Value = 0.853448275862069
$ws.Range("H134").Value = 2.47
$ws.Range("I134").Value = 1.684
$ws.Range("L134").Value = 0.285
$ws.Range("M134").Value = 1.2948
$ws.Range("N134").Value = 1.436
$ws.Range("O134").Value = 1806649.9
$ws.Range("P134").Value = 5525867
$ws.Range("Q134").Value = "Manawatu District"
$ws.Range("R134").Value = "Manawatū"
$ws.Range("S134").Value = "Oroua"
$ws.Range("T134").Value = "Mana_12c"
$ws.Range("U134").Value = "g/m3"

# Row 135
$ws.Range("A135").Value = "Oroua at Mangawhata"
$ws.Range("B135").Value = "Total Nitrogen (Median)"
$ws.Range("D135").Value = "2019 - 2023"
$ws.Range("E135").Value = "RepSite"
$ws.Range("F135").Value = 0.865
$ws.Range("G135").Value = 0.853448275862069
$ws.Range("H135").Value = 2.47
$ws.Range("I135").Value = 1.684
$ws.Range("L135").Value = 0.285
$ws.Range("M135").Value = 1.2948
$ws.Range("N135").Value = 1.436
$ws.Range("O135").Value = 1806649.9
$ws.Range("P135").Value = 5525867
$ws.Range("Q135").Value = "Manawatu District"
$ws.Range("R135").Value = "Manawatū"
$ws.Range("S135").Value = "Oroua"
$ws.Range("T135").Value = "Mana_12c"
$ws.Range("U135").Value = "g/m3"

# Row 136
$ws.Range("A136").Value = "Oroua at Mangawhata"
$ws.Range("B136").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("D136").Value = "2019 - 2023"
$ws.Range("E136").Value = "RepSite"
$ws.Range("F136").Value = 0.054
$ws.Range("G136").Value = 0.118327586206897
$ws.Range("H136").Value = 0.654
$ws.Range("I136").Value = 0.48
$ws.Range("L136").Value = 0.0505
$ws.Range("M136").Value = 0.18356
$ws.Range("N136").Value = 0.37012
$ws.Range("O136").Value = 1806649.9
$ws.Range("P136").Value = 5525867
$ws.Range("Q136").Value = "Manawatu District"
$ws.Range("R136").Value = "Manawatū"
$ws.Range("S136").Value = "Oroua"
$ws.Range("T136").Value = "Mana_12c"
$ws.Range("U136").Value = "g/m3"

# Row 137
$ws.Range("A137").Value = "Oroua at Mangawhata"
$ws.Range("B137").Value = "Total Phosphorus (Median)"
$ws.Range("D137").Value = "2019 - 2023"
$ws.Range("E137").Value = "RepSite"
$ws.Range("F137").Value = 0.054
$ws.Range("G137").Value = 0.118327586206897
$ws.Range("H137").Value = 0.654
$ws.Range("I137").Value = 0.48
$ws.Range("L137").Value = 0.0505
$ws.Range("M137").Value = 0.18356
$ws.Range("N137").Value = 0.37012
$ws.Range("O137").Value = 1806649.9
$ws.Range("P137").Value = 5525867
$ws.Range("Q137").Value = "Manawatu District"
$ws.Range("R137").Value = "Manawatū"
$ws.Range("S137").Value = "Oroua"
$ws.Range("T137").Value = "Mana_12c"
$ws.Range("U137").Value = "g/m3"
